$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$lo = $ws.ListObjects.Item(1)

# --- 1. Add EARNED (SP(3-0-0)) credits of 1.25 for the two rows that
#        previously had no EARNED value (rows 53 & 54 -> Nov/Dec 2023).
$ws.Range("C53").Value = 1.25
$ws.Range("C54").Value = 1.25

# --- 2. Insert a new row before the Jan-2024 entry (currently row 55)
#        so a "2024" year-header row can be added, matching the pattern
#        used for 2020/2021/2022/2023 (rows 10/15/28/42).
$ws.Rows.Item(55).Insert()

# Give the freshly inserted row the same formatting as a normal data row
# (copy from the row directly above, which still has the old formatting).
$ws.Range("A54:K54").Copy()
$ws.Range("A55:K55").PasteSpecial(-4122)

# Column A of a year-header row uses a distinct style (quote-prefixed /
# bold) - copy that single cell's format from the existing "2023" header.
$ws.Range("A42").Copy()
$ws.Range("A55").PasteSpecial(-4122)

# Restore the calculated-column formula that PasteSpecial(formats-only)
# does not carry over.
$ws.Range("G55").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# New year-header row content.
$ws.Range("A55").Value = "'2024"

# --- 3. Populate the row that used to be row 55 (Jan 2024), now shifted
#        down to row 56, with the new leave entry.
$ws.Range("B56").Value = "FL(5-0-0)"
$ws.Range("D56").Value = 5
$ws.Range("K56").Value = "2/23,26,27,28,29/2024"

# --- 4. Grow the table to include the newly inserted row and fix up the
#        calculated-column formula on the new last row (row 133), which
#        PasteSpecial / row-shift left pointing at the wrong cell.
$lo.Resize($ws.Range("A8:K133"))
$ws.Range("G133").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# --- 5. Update the active-cell selection to match where the user ended
#        up after inserting the row (K56, the REMARKS cell of the new
#        Jan-2024 row).
$ws.Range("K56").Select()

$wb.Save()
